$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values
# following a re-run of SGNN dialog act annotation.

$ws.Cells.Item(39, 9).Value = "aa"
$ws.Cells.Item(39, 10).Value = "Agree/Accept"
$ws.Cells.Item(43, 9).Value = "sd"
$ws.Cells.Item(43, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(63, 9).Value = "sd"
$ws.Cells.Item(63, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(64, 9).Value = "sd"
$ws.Cells.Item(64, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(73, 9).Value = "aa"
$ws.Cells.Item(73, 10).Value = "Agree/Accept"
$ws.Cells.Item(86, 9).Value = "aa"
$ws.Cells.Item(86, 10).Value = "Agree/Accept"
$ws.Cells.Item(104, 9).Value = "sd"
$ws.Cells.Item(104, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(119, 9).Value = "sd"
$ws.Cells.Item(119, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(120, 9).Value = "sv"
$ws.Cells.Item(120, 10).Value = "Statement-opinion"
$ws.Cells.Item(123, 9).Value = "sd"
$ws.Cells.Item(123, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(125, 9).Value = "ba"
$ws.Cells.Item(125, 10).Value = "Appreciation"
$ws.Cells.Item(136, 9).Value = "sv"
$ws.Cells.Item(136, 10).Value = "Statement-opinion"
$ws.Cells.Item(154, 9).Value = "b"
$ws.Cells.Item(154, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(158, 9).Value = "sd"
$ws.Cells.Item(158, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(159, 9).Value = "aa"
$ws.Cells.Item(159, 10).Value = "Agree/Accept"
$ws.Cells.Item(160, 9).Value = "aa"
$ws.Cells.Item(160, 10).Value = "Agree/Accept"
$ws.Cells.Item(163, 9).Value = "sv"
$ws.Cells.Item(163, 10).Value = "Statement-opinion"
$ws.Cells.Item(167, 9).Value = "sd"
$ws.Cells.Item(167, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(172, 9).Value = "sd"
$ws.Cells.Item(172, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(183, 9).Value = "sd"
$ws.Cells.Item(183, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(184, 9).Value = "sd"
$ws.Cells.Item(184, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(201, 9).Value = "sv"
$ws.Cells.Item(201, 10).Value = "Statement-opinion"
$ws.Cells.Item(207, 9).Value = "sv"
$ws.Cells.Item(207, 10).Value = "Statement-opinion"
$ws.Cells.Item(208, 9).Value = "%"
$ws.Cells.Item(208, 10).Value = "Uninterpretable"
$ws.Cells.Item(215, 9).Value = "sv"
$ws.Cells.Item(215, 10).Value = "Statement-opinion"
$ws.Cells.Item(216, 9).Value = "sv"
$ws.Cells.Item(216, 10).Value = "Statement-opinion"
$ws.Cells.Item(219, 9).Value = "ba"
$ws.Cells.Item(219, 10).Value = "Appreciation"
$ws.Cells.Item(244, 9).Value = "aa"
$ws.Cells.Item(244, 10).Value = "Agree/Accept"
$ws.Cells.Item(248, 9).Value = "ba"
$ws.Cells.Item(248, 10).Value = "Appreciation"
$ws.Cells.Item(252, 9).Value = "aa"
$ws.Cells.Item(252, 10).Value = "Agree/Accept"
$ws.Cells.Item(253, 9).Value = "b"
$ws.Cells.Item(253, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(255, 9).Value = "aa"
$ws.Cells.Item(255, 10).Value = "Agree/Accept"
$ws.Cells.Item(277, 9).Value = "b"
$ws.Cells.Item(277, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(289, 9).Value = "sd"
$ws.Cells.Item(289, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(290, 9).Value = "sv"
$ws.Cells.Item(290, 10).Value = "Statement-opinion"
$ws.Cells.Item(298, 9).Value = "ba"
$ws.Cells.Item(298, 10).Value = "Appreciation"
$ws.Cells.Item(299, 9).Value = "aa"
$ws.Cells.Item(299, 10).Value = "Agree/Accept"
$ws.Cells.Item(314, 9).Value = "sv"
$ws.Cells.Item(314, 10).Value = "Statement-opinion"
$ws.Cells.Item(317, 9).Value = "sd"
$ws.Cells.Item(317, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(336, 9).Value = "sd"
$ws.Cells.Item(336, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(340, 9).Value = "ba"
$ws.Cells.Item(340, 10).Value = "Appreciation"
$ws.Cells.Item(349, 9).Value = "sd"
$ws.Cells.Item(349, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(358, 9).Value = "b"
$ws.Cells.Item(358, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(365, 9).Value = "sd"
$ws.Cells.Item(365, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(370, 9).Value = "aa"
$ws.Cells.Item(370, 10).Value = "Agree/Accept"
$ws.Cells.Item(372, 9).Value = "ba"
$ws.Cells.Item(372, 10).Value = "Appreciation"
$ws.Cells.Item(386, 9).Value = "sv"
$ws.Cells.Item(386, 10).Value = "Statement-opinion"
$ws.Cells.Item(420, 9).Value = "sv"
$ws.Cells.Item(420, 10).Value = "Statement-opinion"
$ws.Cells.Item(434, 9).Value = "sd"
$ws.Cells.Item(434, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(435, 9).Value = "ba"
$ws.Cells.Item(435, 10).Value = "Appreciation"
$ws.Cells.Item(443, 9).Value = "sd"
$ws.Cells.Item(443, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(446, 9).Value = "sd"
$ws.Cells.Item(446, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(463, 9).Value = "sd"
$ws.Cells.Item(463, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(464, 9).Value = "sd"
$ws.Cells.Item(464, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(466, 9).Value = "sd"
$ws.Cells.Item(466, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(470, 9).Value = "ba"
$ws.Cells.Item(470, 10).Value = "Appreciation"
$ws.Cells.Item(476, 9).Value = "sd"
$ws.Cells.Item(476, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(484, 9).Value = "sd"
$ws.Cells.Item(484, 10).Value = "Statement-non-opinion"
